$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 182, shifting existing rows 182:268 down to 183:269.
$ws.Rows("182:182").Insert()

# Populate the newly inserted row 182 with the new weekly record.
$ws.Range("A182").Value = 3
$ws.Range("B182").Value = "Femacal de La Calera"
$ws.Range("C182").Value = "Coquimbo"
$ws.Range("D182").Value = 44609
$ws.Range("E182").Value = 5
$ws.Range("F182").Value = 100112039
$ws.Range("G182").Value = "Ciboulette"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 120
$ws.Range("K182").Value = 1500
$ws.Range("L182").Value = 1500
$ws.Range("M182").Value = 1500
$ws.Range("N182").Value = "`$/docena de atados"
$ws.Range("O182").Value = "Provincia de Quillota"
$ws.Range("P182").Value = 500
$ws.Range("Q182").Value = 3
$ws.Range("R182").Value = "Hortaliza"
